# "Upload new version with timestamp"
# Fills in the first (and only) sale line item of the day-sale report in
# row 7, the rolled-up selling-price total in row 8, and refreshes the
# generated-at timestamp (plus the page-indicator / footer cells that sit
# next to it) in row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: sale line item ---
$ws.Range("A7").Value = 1

# Name / balance / price / transactions columns hold text, some of which
# looks numeric ("204.00", "1:1", "0:1") - switch those to Text format so
# the values are preserved verbatim as strings instead of being coerced to
# numbers.
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("Q7").NumberFormat = "@"

$ws.Range("C7").Value = "MELLITOFIX 10MG 30 F.C. TABS"   # item name
$ws.Range("H7").Value = "1:1"                             # current balance
$ws.Range("N7").Value = "204.00"                           # price
$ws.Range("Q7").Value = "0:1"                               # transactions

# Order limit / selling price keep their original number formats but are
# still written as text (quote-prefixed) so "1" / "67.3200" round-trip as
# literal strings rather than being normalised to 1 / 67.32.
$ws.Range("L7").Value = "'1"
$ws.Range("P7").Value = "'67.3200"

# --- Row 8: numeric selling-price total (merged P8:Q8, anchor P8) ---
$ws.Range("P8").Value = 67.319999999999993

# --- Row 9: refreshed timestamp / page indicator / footer ---
$ws.Range("A9").Value = "Thursday, 18 September, 2025 9:43 AM"
$ws.Range("G9").Value = "1/1"
$ws.Range("K9").Value = "developed by : Abdelaziz Talaat"
